$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weights")

$ws.Range("B4").Value = -0.1937348010753318
$ws.Range("C4").Value = -6.785386184910475
$ws.Range("D4").Value = 1.689551354844047
$ws.Range("E4").Value = -6.297125354620201
$ws.Range("F4").Value = 7.403454951011472
$ws.Range("G4").Value = -12.85121242792387

$ws.Range("B5").Value = 4.990699365056803
$ws.Range("C5").Value = -8.508136534546692
$ws.Range("D5").Value = -2.207094393071722
$ws.Range("E5").Value = 5.407402955309999
$ws.Range("F5").Value = 10.36228188534344
$ws.Range("G5").Value = -12.99689402641534

$ws.Range("B6").Value = -0.4726016765291365
$ws.Range("C6").Value = -4.016155106674934
$ws.Range("D6").Value = 1.427217109486305
$ws.Range("E6").Value = -0.722224203175398
$ws.Range("F6").Value = 9.189570846447866
$ws.Range("G6").Value = -9.780724066055841

$ws.Range("F7").Value = -2.741007145283092
$ws.Range("G7").Value = 6.18552779297814

$ws.Range("F8").Value = 3.04798120154035
$ws.Range("G8").Value = -3.785340261565299
